$d = $word.ActiveDocument

# Locate the last entry currently in the log: "Finished styling on RSA ..."
$anchorText = "Finished styling on RSA and made it follow with the design and made it responsive. Small polishing changes."
$r = $d.Content
$found = $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found text (still inside the original paragraph,
    # before its trailing bookmark) and split off a brand-new paragraph after it.
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    # The new paragraph's body starts right after the paragraph mark we just added.
    $p1Start = $r.End + 1
    $p1 = $d.Range($p1Start, $p1Start)
    $p1.InsertAfter("17.12.2018 9.45-14.00")

    # Split off a second new paragraph after the first one for the work entry.
    $p1.Collapse(0)
    $p1.InsertParagraphAfter()

    $p2Start = $p1.End + 1
    $p2 = $d.Range($p2Start, $p2Start)
    $p2.InsertAfter("Uploaded the truth table with input to GitHub (not working properly currently) and made sure everything is working correctly and minor polishing changes.")
}
